$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price strings so values like
# "8.630" keep their trailing zero instead of becoming the number 8.63.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values.
$ws.Range("D2").Value = '28.451.60'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.798.41'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '316.69'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").Value = '0.5423'
$ws.Range("E7").Value = '  -1.60%  '
$ws.Range("D8").Value = '0.3778'
$ws.Range("E8").Value = '  -1.88%  '
$ws.Range("D9").Value = '0.07487'
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("D10").Value = '41.94'
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").Value = '20.67'
$ws.Range("E13").Value = '  -2.62%  '
$ws.Range("D14").Value = '6.156'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").Value = '7.301'
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '1.792.89'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '89.54'
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").Value = '0.06518'
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("D22").Value = '5.944'
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("D23").Value = '28.466.37'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '11.11'
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("D25").Value = '2.079'
$ws.Range("E25").Value = '  -1.72%  '
$ws.Range("D26").Value = '159.33'
$ws.Range("E26").Value = '  +1.76%  '
$ws.Range("D27").Value = '20.46'
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("D28").Value = '1.999.86'
$ws.Range("E28").Value = '  -0.75%  '
$ws.Range("D29").Value = '2.327'
$ws.Range("E29").Value = '  -4.54%  '
$ws.Range("D30").Value = '122.87'
$ws.Range("E30").Value = '  -0.52%  '
$ws.Range("E31").Value = '  -5.16%  '
$ws.Range("D32").Value = '0.1057'
$ws.Range("E32").Value = '  +2.43%  '
$ws.Range("D33").Value = '5.618'
$ws.Range("E33").Value = '  -2.10%  '
$ws.Range("D34").Value = '3.653'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("B35").Value = 'Algorand'
$ws.Range("C35").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D35").Value = '0.2267'
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.06504'
$ws.Range("E36").Value = '  +4.13%  '
$ws.Range("D37").Value = '0.02297'
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("D38").Value = '8.630'
$ws.Range("E38").Value = '  -3.52%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("E40").Value = '  -3.13%  '
$ws.Range("D41").Value = '0.6211'
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("D42").Value = '1.451'
$ws.Range("E42").Value = '  +4.84%  '
$ws.Range("D43").Value = '1.192'
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").Value = '13.28'
$ws.Range("E45").Value = '  -1.64%  '
$ws.Range("D46").Value = '3.685'
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").Value = '0.5832'
$ws.Range("E47").Value = '  -2.71%  '
$ws.Range("D48").Value = '126.93'
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("D49").Value = '1.209'
$ws.Range("E49").Value = '  +5.35%  '
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").Value = '0.06893'
$ws.Range("E51").Value = '  -0.50%  '
